# "fix order, send sms" - add Code/Time1/Time2 columns to the receivers
# template so the add_receivers upload can also carry an SMS code + two
# send-time windows, and restyle the header font to Arial.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells: D1=Code, E1=Time1, F1=Time2 (new shared strings 6,7,8)
$ws.Range("D1").Value = "Code"
$ws.Range("E1").Value = "Time1"
$ws.Range("F1").Value = "Time2"

# Match the bold header styling already used by A1:C1
$ws.Range("D1:F1").Font.Bold = $true

# Re-brand the sheet font from Calibri to Arial (headers + existing data
# row), without touching the still-empty D2:F2 cells.
$ws.Range("A1:C2").Font.Name = "Arial"
$ws.Range("D1:F1").Font.Name = "Arial"

# Nudge the two data columns a touch narrower to match the refreshed layout.
$ws.Columns.Item(2).ColumnWidth = 14.33
$ws.Columns.Item(3).ColumnWidth = 15.33

# Leave the selection on the new Time1 cell, like the source workbook.
$ws.Range("E2").Select()

$wb.Save()
